# Applies the reordering of the "2928" block (2 rows) and the "2957" block
# (4 rows) within rows 51-56: the 2928 block now comes first, followed by
# the 2957 block. Row contents (B..H) travel together with their row; only
# their position (and therefore row-relative B/group index) changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for rows 51-56, columns A-H, before touching them.
$original = @{}
for ($r = 51; $r -le 56; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le 8; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $rowVals
}

# New order: the block that used to live in rows 55-56 now occupies
# rows 51-52, and the block that used to live in rows 51-54 now occupies
# rows 53-56.
$sourceOrder = @(55, 56, 51, 52, 53, 54)

$destRow = 51
foreach ($srcRow in $sourceOrder) {
    $srcVals = $original[$srcRow]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
    $destRow++
}
